# Level design pass ("Level2"): rebalance the cost multipliers on Sheet3.
# Sheet3!B2:B4 (C*E, shared formula) and Sheet1!E3:E14 (VLOOKUP into
# Sheet3!B2:B4) recalculate automatically from these input changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4.5
$ws.Range("C4").Value = 6

$ws.Range("C5").Select()

$excel.Calculate()
